# Weekly update: insert a new week's price row above the current row 4,
# pushing the existing rows 4-19 down to 5-20.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("4:4").Insert()

# Populate the newly inserted row 4 with this week's data.
$ws.Cells.Item(4, 1).Value  = 7
$ws.Cells.Item(4, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(4, 3).Value  = "Ñuble"
$ws.Cells.Item(4, 4).Value  = (Get-Date -Year 2022 -Month 9 -Day 7 -Hour 0 -Minute 0 -Second 0).Date
$ws.Cells.Item(4, 5).Value  = 16
$ws.Cells.Item(4, 6).Value  = 100112044
$ws.Cells.Item(4, 7).Value  = "Perejil"
$ws.Cells.Item(4, 8).Value  = "Sin especificar"
$ws.Cells.Item(4, 9).Value  = "Primera"
$ws.Cells.Item(4, 10).Value = 300
$ws.Cells.Item(4, 11).Value = 750
$ws.Cells.Item(4, 12).Value = 850
$ws.Cells.Item(4, 13).Value = 800
$ws.Cells.Item(4, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(4, 15).Value = "Región del Maule"
$ws.Cells.Item(4, 16).Value = 800
$ws.Cells.Item(4, 17).Value = 1
$ws.Cells.Item(4, 18).Value = "Hortaliza"
